$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.933
$ws.Range("E3").Value = 16.263
$ws.Range("B21").Value = 9.238
$ws.Range("B23").Value = 7.179
$ws.Range("E24").Value = 16.479
$ws.Range("B25").Value = 6.425999999999999
$ws.Range("D27").Value = -8.517000000000001
$ws.Range("D31").Value = -8.248000000000001
$ws.Range("D39").Value = -8.063000000000001
$ws.Range("D48").Value = -7.31
$ws.Range("D51").Value = -8.373999999999999
$ws.Range("D52").Value = -8.054
$ws.Range("B53").Value = 6.086
$ws.Range("D55").Value = -8.035
$ws.Range("D56").Value = -8.416999999999998
$ws.Range("B57").Value = 4.924999999999999
$ws.Range("D57").Value = -8.059999999999999
$ws.Range("E57").Value = 16.623
$ws.Range("B59").Value = 5.145
$ws.Range("E61").Value = 16.706
$ws.Range("B69").Value = 5.388
$ws.Range("E70").Value = 17.568
$ws.Range("D73").Value = -8.278000000000002
$ws.Range("B79").Value = 5.646000000000001
$ws.Range("B83").Value = 5.915
$ws.Range("E86").Value = 16.597
$ws.Range("D89").Value = -6.173999999999999
$ws.Range("D90").Value = -7.49
$ws.Range("B93").Value = 5.692
$ws.Range("E98").Value = 16.493
$ws.Range("E100").Value = 16.725
$ws.Range("E102").Value = 16.698
